$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data row that referenced the ".COM 2" table type (row 3),
# leaving the formatting/styles intact, as part of the database refresh.
$ws.Range("A3:G3").ClearContents()

# Update the active selection to reflect the new focus cell.
$ws.Range("G3").Select()
